# Database Updated, Update Database.py added to menu
#
# The token data for each card used to be stored one fact per row
# (name, type line, ability/keyword lines, power/toughness all on their
# own rows). The generator script now emits one row per token as a
# Python tuple literal: (name, [type_line, ...other lines..., pt]).
# Collapse the old per-fact rows A2:A29 into the new consolidated rows
# A2:A7, and drop the now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Angel', ['Token Creature — Angel', 'Flying', '4/4'])"
$ws.Range("A3").Value = "('Beast', ['Token Creature — Beast', '4/4'])"
$ws.Range("A4").Value = "('Bird', ['Token Creature — Bird', 'Flying', '1/1'])"
$ws.Range("A5").Value = "('Djinn Monk', ['Token Creature — Djinn Monk', 'Flying', '2/2'])"
$ws.Range("A6").Value = "('Dragon', ['Token Creature — Dragon', 'Flying', '5/5'])"
$ws.Range("A7").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"

# The rest of the old data (rows 8-29) is now obsolete; remove it so the
# sheet's used range shrinks back down to A1:A7.
$ws.Rows("8:29").Delete()
